# "Add files via upload" — re-upload of the "LOS Galacticos" yahoo roster
# sheet with an updated player lineup (players re-ordered / swapped and
# Amen Thompson's eligible positions updated to include PF).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final desired state for A2:C19 (Oyuncu Adı / Pozisyon / Takım).
$rows = @(
    @(2,  "James Harden",          "PG,SG",    "LA Clippers"),
    @(3,  "Bradley Beal",          "PG,SG,SF", "Phoenix Suns"),
    @(4,  "Fred VanVleet",         "PG",       "Houston Rockets"),
    @(5,  "Anthony Edwards",       "SG,SF",    "Minnesota Timberwolves"),
    @(6,  "Amen Thompson",         "SG,SF,PF", "Houston Rockets"),
    @(7,  "Jaren Jackson Jr.",     "PF,C",     "Memphis Grizzlies"),
    @(8,  "Jayson Tatum",          "SF,PF",    "Boston Celtics"),
    @(9,  "P.J. Washington",       "PF",       "Dallas Mavericks"),
    @(10, "Giannis Antetokounmpo", "PF,C",     "Milwaukee Bucks"),
    @(11, "Zion Williamson",       "PF,C",     "New Orleans Pelicans"),
    @(12, "Ivica Zubac",           "C",        "LA Clippers"),
    @(13, "Keyonte George",        "PG,SG",    "Utah Jazz"),
    @(14, "Donovan Clingan",       "C",        "Portland Trail Blazers"),
    @(15, "Robert Williams III",   "C",        "Portland Trail Blazers"),
    @(16, "Anfernee Simons",       "PG,SG",    "Portland Trail Blazers"),
    @(17, "Jonathan Kuminga",      "SF,PF",    "Golden State Warriors"),
    @(18, "Goga Bitadze",          "C",        "Orlando Magic"),
    @(19, "Paul George",           "SG,SF,PF", "Philadelphia 76ers")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}
